# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" (fund holdings detail) right after the
# "总计" (Total) summary sheet, pushing the existing "2021-Q4" / "2021-Q3"
# sheets one slot to the right, and adds a matching summary row on the
# "总计" sheet.

$wb = $excel.ActiveWorkbook

# Keep stable references to the sheets that already exist, captured before
# any insertion so they remain valid no matter how the tab order/index
# shifts afterwards.
$wsTotal = $wb.Worksheets.Item(1)      # "总计"
$wsQ4    = $wb.Worksheets.Item(2)      # "2021-Q4" (style/template donor)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wsTotal)
$wsNew.Name = "2022-Q3"

# Match the page margins used by the other detail sheets.
$wsNew.PageSetup.LeftMargin = 54
$wsNew.PageSetup.RightMargin = 54
$wsNew.PageSetup.TopMargin = 72
$wsNew.PageSetup.BottomMargin = 72
$wsNew.PageSetup.HeaderMargin = 36
$wsNew.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Header row (B1:H1).
# ---------------------------------------------------------------------
$headers = New-Object 'object[,]' 1,7
$headers[0,0] = "基金代码"
$headers[0,1] = "基金名称"
$headers[0,2] = "基金规模"
$headers[0,3] = "股票总仓位"
$headers[0,4] = "仓位占比"
$headers[0,5] = "持有市值(亿元)"
$headers[0,6] = "仓位排名"
$wsNew.Range("B1:H1").Value = $headers

# Copy the bold/centered/bordered header style used throughout the
# workbook (style index carried via the existing "2021-Q4" sheet header).
$wsQ4.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 3. Data rows (A2:H14).
# ---------------------------------------------------------------------
$data = New-Object 'object[,]' 13,8
$data[0,0] = 0
$data[0,1] = "'100061"
$data[0,2] = "富国中国中小盘混合（QDII）人民币"
$data[0,3] = "'35.11"
$data[0,4] = "'83.32"
$data[0,5] = "'2.12"
$data[0,6] = "'0.7443"
$data[0,7] = 10
$data[1,0] = 1
$data[1,1] = "'010591"
$data[1,2] = "富国中国中小盘混合（QDII）美元"
$data[1,3] = "'35.11"
$data[1,4] = "'83.32"
$data[1,5] = "'2.12"
$data[1,6] = "'0.7443"
$data[1,7] = 10
$data[2,0] = 2
$data[2,1] = "'005583"
$data[2,2] = "易方达港股通红利灵活配置混合"
$data[2,3] = "'6.98"
$data[2,4] = "'90.31"
$data[2,5] = "'4.07"
$data[2,6] = "'0.2841"
$data[2,7] = 10
$data[3,0] = 3
$data[3,1] = "'012227"
$data[3,2] = "景顺长城港股通全球竞争力混合A"
$data[3,3] = "'7.81"
$data[3,4] = "'74.96"
$data[3,5] = "'3.51"
$data[3,6] = "'0.2741"
$data[3,7] = 8
$data[4,0] = 4
$data[4,1] = "'160526"
$data[4,2] = "博时优势企业灵活配置混合（LOF）A"
$data[4,3] = "'5.29"
$data[4,4] = "'85.98"
$data[4,5] = "'3.51"
$data[4,6] = "'0.1857"
$data[4,7] = 10
$data[5,0] = 5
$data[5,1] = "'040018"
$data[5,2] = "华安香港精选股票（QDII）"
$data[5,3] = "'4.36"
$data[5,4] = "'85.97"
$data[5,5] = "'3.49"
$data[5,6] = "'0.1522"
$data[5,7] = 5
$data[6,0] = 6
$data[6,1] = "'012228"
$data[6,2] = "景顺长城港股通全球竞争力混合C"
$data[6,3] = "'0.99"
$data[6,4] = "'74.96"
$data[6,5] = "'3.51"
$data[6,6] = "'0.0347"
$data[6,7] = 8
$data[7,0] = 7
$data[7,1] = "'005701"
$data[7,2] = "上投摩根香港精选港股通混合A"
$data[7,3] = "'0.42"
$data[7,4] = "'82.85"
$data[7,5] = "'3.64"
$data[7,6] = "'0.0153"
$data[7,7] = 5
$data[8,0] = 8
$data[8,1] = "'005143"
$data[8,2] = "中融沪港深大消费主题灵活配置混合C"
$data[8,3] = "'0.27"
$data[8,4] = "'90.10"
$data[8,5] = "'4.45"
$data[8,6] = "'0.0120"
$data[8,7] = 8
$data[9,0] = 9
$data[9,1] = "'040021"
$data[9,2] = "华安大中华升级股票（QDII）"
$data[9,3] = "'0.29"
$data[9,4] = "'68.57"
$data[9,5] = "'3.40"
$data[9,6] = "'0.0099"
$data[9,7] = 2
$data[10,0] = 10
$data[10,1] = "'005142"
$data[10,2] = "中融沪港深大消费主题灵活配置混合A"
$data[10,3] = "'0.13"
$data[10,4] = "'90.10"
$data[10,5] = "'4.45"
$data[10,6] = "'0.0058"
$data[10,7] = 8
$data[11,0] = 11
$data[11,1] = "'007234"
$data[11,2] = "博时优势企业灵活配置混合（LOF）C"
$data[11,3] = "'0.00"
$data[11,4] = "'85.98"
$data[11,5] = "'3.51"
$data[11,6] = 0
$data[11,7] = 10
$data[12,0] = 12
$data[12,1] = "'016921"
$data[12,2] = "上投摩根香港精选港股通混合C"
$data[12,3] = "'0.00"
$data[12,4] = "'82.85"
$data[12,5] = "'3.64"
$data[12,6] = 0
$data[12,7] = 5

$wsNew.Range("A2:H14").Value = $data

# Column A (row index numbers) shares the same bold/centered/bordered
# style as the header cells and as the "总计" index column.
$wsQ4.Range("A2").Copy()
$wsNew.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 4. Update the "总计" (Total) summary sheet: insert a new row for the
#    2022-Q3 figures above the existing 2021-Q4 / 2021-Q3 rows.
# ---------------------------------------------------------------------
$wsTotal.Rows(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)    # xlPasteFormats

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 13
$wsTotal.Range("D2").Value = 2.46

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

Write-Output "2022-Q3 sheet added and 总计 sheet updated."
